$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace every occurrence of "DH3" with "DH1" across the table's text cells.
$ws.Range("C2").Value  = "DH1项目主计划"
$ws.Range("F2").Value  = "【飞书】DH1项目主计划`n12月8日版本"
$ws.Range("C3").Value  = "DH1整车配置清单"
$ws.Range("F3").Value  = "【飞书】DH1整车配置表`n12月3日版本"
$ws.Range("F4").Value  = "【飞书】包含DH1，DH1各业务板块联系人信息"
$ws.Range("C5").Value  = "DH1电子电气拓扑图"
$ws.Range("C6").Value  = "DH1项目J2阶段评审材料输入信息"
$ws.Range("D6").Value  = '\\10.4.9.25\Project\DH系列\060-整车开发\DH1\J2评审输入'
$ws.Range("C7").Value  = "DH1项目J2阶段评审材料归档路径"
$ws.Range("D7").Value  = '\\10.4.9.25\Project\DH系列\130-智能软件\DH1\130100_ 阶段评审材料\130102_ J2阶段评审材料及会议纪要'
$ws.Range("C8").Value  = "DH1项目智软控制器开发计划"
$ws.Range("C9").Value  = "DH1项目华为接口联系人清单"
$ws.Range("C10").Value = "DH1项目研发总院智软内部责任分工"
$ws.Range("C11").Value = "DH1不满再发防止-整车级"
$ws.Range("F11").Value = "【飞书】DH1不满再发防止-整车级"
$ws.Range("C12").Value = "DH1 LLR和PRC再发防止 - 总院级"
$ws.Range("F12").Value = "【企微】DH1 LLR和PRC再发防止 - 总院级"

# Re-fit rows whose text now wraps onto two lines again, so the row
# height stays at its original (non-custom) state.
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()

# Update the last active-cell selection recorded in the sheet view.
$ws.Range("D20").Select()
